$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.023021697998047
$ws.Range("B1").Value = 3.392872095108032
$ws.Range("C1").Value = 3.028100967407227
$ws.Range("D1").Value = 3.27067232131958
$ws.Range("E1").Value = -1
